# Swap the contents of row 2 and row 4 in the "Artfynd" sheet.
# (The two occurrence records were re-ordered; row 3 is untouched.)
# Columns that change: A, B, D, E, F, G, H, Q, R, and the AC "public
# comment" field, which only exists on the record that ends up in row 2
# in the "before" state / row 4 in the "after" state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values currently in row 2 (to be moved to row 4) ---
$A2 = 111399971
$B2 = 90666
$D2 = "LC"
$E2 = 4364
$F2 = "Dropptaggsvamp"
$G2 = "Hydnellum ferrugineum"
$H2 = "(Fr.:Fr.) P. Karst."
$Q2 = 713904.7870529388
$R2 = 7286390.269694135
$AC2 = "fjolårets fruktkropp, vit mycelmatta i marken"

# --- values currently in row 4 (to be moved to row 2) ---
$A4 = 111399984
$B4 = 89558
$D4 = "VU"
$E4 = 1503
$F4 = "Gräddporing"
$G4 = "Sidera lenis"
$H4 = "(P.Karst.) Miettinen"
$Q4 = 713925.3805777475
$R4 = 7286420.397548387

# --- write row 2 <= old row 4 values ---
$ws.Range("A2").Value = $A4
$ws.Range("B2").Value = $B4
$ws.Range("D2").Value = $D4
$ws.Range("E2").Value = $E4
$ws.Range("F2").Value = $F4
$ws.Range("G2").Value = $G4
$ws.Range("H2").Value = $H4
$ws.Range("Q2").Value = $Q4
$ws.Range("R2").Value = $R4
$ws.Range("AC2").Value = ""

# --- write row 4 <= old row 2 values ---
$ws.Range("A4").Value = $A2
$ws.Range("B4").Value = $B2
$ws.Range("D4").Value = $D2
$ws.Range("E4").Value = $E2
$ws.Range("F4").Value = $F2
$ws.Range("G4").Value = $G2
$ws.Range("H4").Value = $H2
$ws.Range("Q4").Value = $Q2
$ws.Range("R4").Value = $R2
$ws.Range("AC4").Value = $AC2
